# Pavel - new user for linking test
# Adds a new "Linking_AutoUser" row to the Users sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$newRow = 52

$ws.Cells.Item($newRow, 1).Value = "Linking_AutoUser"
$ws.Cells.Item($newRow, 2).Value = "Password1"
$ws.Cells.Item($newRow, 5).Value = "Default user for Linking tests"
$ws.Cells.Item($newRow, 6).Value = "N"
$ws.Cells.Item($newRow, 7).Value = "linking.autouser@mailinator.com"

# Apply the same thin all-round bordered style used throughout the table
# rows (matches the style used on row 29, a non-hyperlinked data row).
$ws.Range("A29:G29").Copy()
$ws.Range("A52:G52").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Scroll the view and select the newly added row, matching the saved view state
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A52:G52").Select()
